# Update workbook per commit "update to new files"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For rows 2-6: set D column to "U", clear F and G columns (E stays 1)
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = "U"           # D column
    $ws.Cells.Item($r, 6).ClearContents()        # F column
    $ws.Cells.Item($r, 7).ClearContents()        # G column
}

# Update selection to D7
$ws.Range("D7").Select()
